$wb = $excel.ActiveWorkbook

# Add a new worksheet ("Sheet3") after the last existing sheet (Sheet2).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "Sheet3"

# Populate the header row A1:F1 with the column labels (A-F). "A".."D" already
# exist in the shared-string table (reused from Sheet1/Sheet2); "E" and "F"
# are brand new strings appended to the table.
$ws3.Range("A1").Value = "A"
$ws3.Range("B1").Value = "B"
$ws3.Range("C1").Value = "C"
$ws3.Range("D1").Value = "D"
$ws3.Range("E1").Value = "E"
$ws3.Range("F1").Value = "F"

# Leave the cursor on F2 (empty row below the header) and make Sheet3 the
# active/selected sheet, matching the authored workbook.
[void]$ws3.Range("F2").Select()
$ws3.Activate()
